$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131
$xlPasteValues = -4163

# Rows where column H currently holds a numeric percentage that mirrors the
# text already stored in column G (e.g. H2 = 1E-3 / 0.1%, G2 = "0.1%").
# These get converted to text (matching G) and left-aligned.
$textRows = 2,4,5,6,7,8,9

foreach ($r in $textRows) {
    $src = $ws.Range("G$r")
    $dst = $ws.Range("H$r")

    # Drop H's existing (percentage) number format so the upcoming paste
    # doesn't inherit it - we want a plain "General" text cell like G.
    $dst.ClearFormats() | Out-Null

    # Copy the text value from G into H (values only, so no formatting is
    # dragged along), then left-align both source and destination cells.
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteValues) | Out-Null

    $src.HorizontalAlignment = $xlLeft
    $dst.HorizontalAlignment = $xlLeft
}

# Row 3 keeps its numeric percentage values in both G and H; only the
# alignment changes to left.
$ws.Range("G3").HorizontalAlignment = $xlLeft
$ws.Range("H3").HorizontalAlignment = $xlLeft

# Update the sheet's last selected cell.
$ws.Range("K9").Select() | Out-Null
